# Generate Report for Handback
# The 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.md file has been handed back and is
# now in sync with en-US (previously it was only "Ready for handoff").
# Update the Overview sheet's status columns and each locale sheet's Status /
# Latest Handback DateTime columns for that file's row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: row for 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.md is row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: row for 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.md is row 3 ---
# Column C = Status, Column H = Latest Handback DateTime
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("H3").Value = "2016-03-18 16:41:21"

# --- de-de sheet: row for 9ed7af9a-55b6-484b-ba97-0453c63e4a9a.md is row 3 ---
# Column C = Status, Column H = Latest Handback DateTime
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("H3").Value = "2016-03-18 16:41:27"
